$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.046.93"
$ws.Range("E2").Value = "  -2.45%  "

$ws.Range("D3").Value = "1.829.73"
$ws.Range("E3").Value = "  -1.28%  "

$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.87%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3862"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.75%  "

$ws.Range("E9").Value = "  -0.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9592"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.87"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.98%  "

$ws.Range("D12").Value = "1.839.04"
$ws.Range("E12").Value = "  -0.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.656"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.884"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06782"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009949"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.85%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("D21").Value = "28.074.29"
$ws.Range("E21").Value = "  -2.37%  "

$ws.Range("E22").Value = "  -2.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.104"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.22%  "

$ws.Range("D25").Value = "2.115.08"
$ws.Range("E25").Value = "  +1.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.37%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.715"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.971"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9357"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.46%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09244"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.299"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.316"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.323"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.95%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05867"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02143"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.144"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.725"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5590"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.899"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1763"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.217"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5267"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07012"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.82%  "

$ws.Range("E47").Value = "  -10.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.834"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.85%  "

$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.320"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.37%  "
